$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("VerifyDeleteOffer")

# Clear the "Status" (column H) results that were previously marked Pass/Fail,
# restoring the cells to their blank state while keeping formatting/style intact.
$ws.Range("H5:H30").ClearContents()
$ws.Range("H32:H35").ClearContents()
